# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 5
    19 = 3
    20 = 2
    21 = 0
    22 = 4
    23 = 1
    24 = 1
    25 = 1
    26 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
